# Incomes sheet: populate the income records table (Source, Amount, Date)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout:
#   A: Source (text)
#   B: Amount (number)
#   C: Date   (text - stored verbatim, e.g. "13/8/2025", not converted to a date serial)

$sources = @("Source", "salary", "toy", "girlfriend", "bakchodi", "freelance", "Avacado", "sold car", "abhayawas@gmail.com", "john@example.com")
$amounts = @("Amount", 2000, 5, 500, 2000, 679, 10, 5000, 1000, 1000)
$dates   = @("Date", "6/9/2025", "13/8/2025", "13/8/2025", "13/8/2025", "12/8/2025", "12/8/2025", "7/7/2025", "3/6/2025", "3/5/2024")

$rowCount = $sources.Length

# Pre-format A and C as Text so the date-looking strings aren't reinterpreted
# as date serial numbers when the value is assigned.
$ws.Range("A1:A$rowCount").NumberFormat = "@"
$ws.Range("C1:C$rowCount").NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $sources[$i]
    $ws.Cells.Item($row, 2).Value = $amounts[$i]
    $ws.Cells.Item($row, 3).Value = $dates[$i]
}
